$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the stray semicolon in A148 so this cell's text matches the existing
# "Ctrl+ up / down" shared string (row 48) instead of its own unique one.
$ws.Range("A148").Value2 = "Ctrl+ up / down"

# Delete column D ("grade"/"grade1"/"grade2") entirely; column E ("grupo")
# shifts left to become the new column D.
$ws.Columns("D").Delete()

# The autofilter button row is turned off (underlying data range shrank).
$ws.AutoFilterMode = $false

# Keep the hidden _FilterDatabase defined name in sync with the new extent.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='store-data'!`$A`$1:`$D`$150"
    }
}

# Restore the selection to the cell the author ended up on after editing.
$null = $ws.Range("D134").Select()
